# Update the Andre Russell activity sheet with the latest values pulled
# from the updated Excel form. Columns: C=runs, D=balls, E=fours, F=sixes.
# Values are stored as text (matching the sheet's existing
# "numbers stored as text" convention), so each value is entered with a
# leading apostrophe to force text entry, then the style is reset to
# "Normal" so no stray quote-prefix formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.Value = "'" + $value
    $rng.Style = "Normal"
}

Set-TextValue "C2" "12"
Set-TextValue "D2" "9"

Set-TextValue "C3" "2"
Set-TextValue "D3" "4"
Set-TextValue "F3" "0"

Set-TextValue "C4" "25"
Set-TextValue "F4" "3"

Set-TextValue "C5" "5"
Set-TextValue "D5" "3"
Set-TextValue "F5" "0"

Set-TextValue "C6" "16"
Set-TextValue "D6" "10"
Set-TextValue "E6" "2"

Set-TextValue "C7" "9"
Set-TextValue "D7" "11"
Set-TextValue "E7" "1"
Set-TextValue "F7" "0"

Set-TextValue "C9" "13"
Set-TextValue "D9" "8"
Set-TextValue "E9" "1"
Set-TextValue "F9" "1"

Set-TextValue "C10" "24"
Set-TextValue "D10" "14"
Set-TextValue "E10" "0"
Set-TextValue "F10" "3"
